$wb = $excel.ActiveWorkbook

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 81.888885
$ws.Range("I9").Value = 72.833336
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 72.833336
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 96.166664
$ws.Range("N9").Value = -438

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1622.8
$ws.Range("I18").Value = 1103.2307
$ws.Range("K18").Value = 1103.2307
$ws.Range("M18").Value = -819.2307000000001

# ALC row 46
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1965.3334
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1965.3334
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5896.0002
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -6134.0002

# ALC row 60
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 1965.3334
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1965.3334
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 5896.0002
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -6864.0002

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 698
$ws.Range("J70").Value = 777.2
$ws.Range("L70").Value = 2331.6
$ws.Range("N70").Value = -2871.6

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 698
$ws.Range("J73").Value = 777.2
$ws.Range("L73").Value = 2331.6
$ws.Range("N73").Value = -4203.6

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3268.2415
$ws.Range("I74").Value = 3245.55
$ws.Range("J74").Value = 3318.6667
$ws.Range("K74").Value = 3245.55
$ws.Range("L74").Value = 3318.6667
$ws.Range("M74").Value = -2309.55
$ws.Range("N74").Value = -5190.6667

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3268.2415
$ws.Range("I77").Value = 3245.55
$ws.Range("J77").Value = 3318.6667
$ws.Range("K77").Value = 16227.75
$ws.Range("L77").Value = 16593.3335
$ws.Range("M77").Value = -11547.75
$ws.Range("N77").Value = -25953.3335

# ALC row 114
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 40722
$ws.Range("J114").Value = 40722
$ws.Range("L114").Value = 40722
$ws.Range("N114").Value = -49400

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1838981.8
$ws.Range("I132").Value = 2119425
$ws.Range("J132").Value = 521.3333
$ws.Range("K132").Value = 6358275
$ws.Range("L132").Value = 1563.9999
$ws.Range("M132").Value = -6355745
$ws.Range("N132").Value = -6623.9999

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5110.625
$ws.Range("I141").Value = 5836.5386
$ws.Range("J141").Value = 1965
$ws.Range("K141").Value = 17509.6158
$ws.Range("L141").Value = 5895
$ws.Range("M141").Value = -12329.6158
$ws.Range("N141").Value = -16255

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2081.975
$ws.Range("I63").Value = 1942.1072
$ws.Range("K63").Value = 1942.1072
$ws.Range("M63").Value = -1256.1072

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2081.975
$ws.Range("I66").Value = 1942.1072
$ws.Range("K66").Value = 9710.536
$ws.Range("M66").Value = -6278.536

# ARM row 103
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# BSM row 6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 11928
$ws.Range("J6").Value = 11928
$ws.Range("L6").Value = 11928
$ws.Range("N6").Value = -12154

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3649.375
$ws.Range("I86").Value = 3628.1316
$ws.Range("J86").Value = 3730.1
$ws.Range("K86").Value = 3628.1316
$ws.Range("L86").Value = 3730.1
$ws.Range("M86").Value = -2505.1316
$ws.Range("N86").Value = -5976.1

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3649.375
$ws.Range("I89").Value = 3628.1316
$ws.Range("J89").Value = 3730.1
$ws.Range("K89").Value = 18140.658
$ws.Range("L89").Value = 18650.5
$ws.Range("M89").Value = -12524.658
$ws.Range("N89").Value = -29882.5

# CRP row 100
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 40791.332
$ws.Range("J100").Value = 40791.332
$ws.Range("L100").Value = 40791.332
$ws.Range("N100").Value = -42955.332

# CUL row 94
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3170.0667
$ws.Range("I94").Value = 1500
$ws.Range("J94").Value = 3289.3572
$ws.Range("K94").Value = 4500
$ws.Range("L94").Value = 9868.071599999999
$ws.Range("M94").Value = -3824
$ws.Range("N94").Value = -11220.0716

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 902.86
$ws.Range("I131").Value = 514.875
$ws.Range("J131").Value = 936.59784
$ws.Range("K131").Value = 1544.625
$ws.Range("L131").Value = 2809.79352
$ws.Range("M131").Value = 3495.375
$ws.Range("N131").Value = -12889.79352

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4869.5
$ws.Range("I133").Value = 3453.077
$ws.Range("J133").Value = 7500
$ws.Range("K133").Value = 10359.231
$ws.Range("L133").Value = 22500
$ws.Range("M133").Value = -5299.231
$ws.Range("N133").Value = -32620

# LTW row 58
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 1062
$ws.Range("I58").Value = 1062
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1062
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -802
$ws.Range("N58").ClearContents()

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2398.1042
$ws.Range("I68").Value = 877.61536
$ws.Range("J68").Value = 2962.8572
$ws.Range("K68").Value = 877.61536
$ws.Range("L68").Value = 2962.8572
$ws.Range("M68").Value = -128.61536
$ws.Range("N68").Value = -4460.8572

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2398.1042
$ws.Range("I71").Value = 877.61536
$ws.Range("J71").Value = 2962.8572
$ws.Range("K71").Value = 4388.0768
$ws.Range("L71").Value = 14814.286
$ws.Range("M71").Value = -644.0767999999998
$ws.Range("N71").Value = -22302.286

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2103.3171
$ws.Range("I100").Value = 1478
$ws.Range("J100").Value = 2362.0688
$ws.Range("K100").Value = 1478
$ws.Range("L100").Value = 2362.0688
$ws.Range("M100").Value = -937
$ws.Range("N100").Value = -3444.0688

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3807.8125
$ws.Range("I122").Value = 4118.75
$ws.Range("K122").Value = 12356.25
$ws.Range("M122").Value = -9906.25

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 26501.562
$ws.Range("J133").Value = 26501.562
$ws.Range("L133").Value = 26501.562
$ws.Range("N133").Value = -31561.562

# WVR row 27
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 38313
$ws.Range("J27").Value = 38313
$ws.Range("L27").Value = 38313
$ws.Range("N27").Value = -38451

# WVR row 125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 40650
$ws.Range("J125").Value = 40650
$ws.Range("L125").Value = 40650
$ws.Range("N125").Value = -50490
